# Update the base metric values on the "Metrics" sheet. These feed the
# formulas on the "today" sheet (B11:B22, and their E/F mirrors), so those
# will recalculate automatically once these are written.
$wb = $excel.ActiveWorkbook
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 375280.41000000003
$metrics.Range("B3").Value = 282631.43000000005
$metrics.Range("B4").Value = 97627.57
$metrics.Range("B5").Value = 15323
$metrics.Range("B6").Value = 6011151.1400000006
$metrics.Range("B7").Value = 5053349.0599999996
$metrics.Range("B8").Value = 1761719.39
$metrics.Range("B9").Value = 235600
$metrics.Range("B10").Value = 34476532.129999995
$metrics.Range("B11").Value = 32328624.219999999
$metrics.Range("B12").Value = 12043441.43
$metrics.Range("B13").Value = 1333230

# Move the selection on "Metrics" (previously the active tab) off its old
# cell, then switch the active tab to "today" and move its selection too -
# matches the workbookView/sheetView changes in the diff (tabSelected moves
# from Metrics to today, activeTab becomes 5).
$metrics.Activate()
$metrics.Range("D20").Select()

$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("E7").Select()
